$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.667.01"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "3.494.75"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("E7").Value = "  -2.19%  "
$ws.Range("D8").Value = "3.486.76"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.205"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.56%  "
$ws.Range("E11").Value = "  -3.99%  "
$ws.Range("E12").Value = "  -4.39%  "
$ws.Range("E13").Value = "  -6.35%  "
$ws.Range("E14").Value = "  -4.22%  "
$ws.Range("D15").Value = "4.046.79"
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "648.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.29%  "
$ws.Range("D17").Value = "69.722.05"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "3.497.56"
$ws.Range("E18").Value = "  -3.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.40%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.949"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.94%  "
$ws.Range("E26").Value = "  -7.30%  "
$ws.Range("E27").Value = "  -4.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.42%  "
$ws.Range("E31").Value = "  -6.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.46%  "
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("E34").Value = "  -5.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "61.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "548.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.72%  "
$ws.Range("D37").Value = "3.698.21"
$ws.Range("E37").Value = "  -6.50%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "0.0₃0789"
$ws.Range("E39").Value = "  -9.11%  "
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.372"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +39.89%  "
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.32%  "
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.21%  "
$ws.Range("E49").Value = "  -4.15%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.06%  "
